# Auto-generated Excel COM-interop script applying the diff to Seraph_Profits workbook.
# Each block targets one (sheet, row) pair; per-cell operations either set a new value
# or clear the cell entirely (for cells removed in the diff) / add a new cell (for cells
# newly introduced in the diff).

$wb = $excel.ActiveWorkbook

# Hunk 0: sheet ALC, row 47
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

# Hunk 1: sheet ALC, row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3200
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

# Hunk 2: sheet ALC, row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3200
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

# Hunk 3: sheet ALC, row 81
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 37861.5
$ws.Range("J81").Value = 37861.5
$ws.Range("L81").Value = 37861.5
$ws.Range("N81").Value = -39857.5

# Hunk 4: sheet ALC, row 84
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H84").Value = 37861.5
$ws.Range("J84").Value = 37861.5
$ws.Range("L84").Value = 113584.5
$ws.Range("N84").Value = -123568.5

# Hunk 5: sheet ALC, row 97
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 4951.9
$ws.Range("J97").Value = 4566.5557
$ws.Range("L97").Value = 13699.6671
$ws.Range("N97").Value = -14691.6671

# Hunk 6: sheet ALC, row 118
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()

# Hunk 7: sheet ALC, row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1032.238
$ws.Range("I135").Value = 812.3333
$ws.Range("J135").Value = 1582
$ws.Range("K135").Value = 7310.9997
$ws.Range("L135").Value = 14238
$ws.Range("M135").Value = -4775.9997
$ws.Range("N135").Value = -19308

# Hunk 8: sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4885.857
$ws.Range("I137").Value = 2127.0908
$ws.Range("J137").Value = 7920.5
$ws.Range("K137").Value = 6381.2724
$ws.Range("L137").Value = 23761.5
$ws.Range("M137").Value = -3831.2724
$ws.Range("N137").Value = -28861.5

# Hunk 9: sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6330.7646
$ws.Range("J138").Value = 5749.3125
$ws.Range("L138").Value = 17247.9375
$ws.Range("N138").Value = -27527.9375

# Hunk 10: sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17772.791
$ws.Range("I32").Value = 8370.322
$ws.Range("J32").Value = 27175.258
$ws.Range("K32").Value = 8370.322
$ws.Range("L32").Value = 27175.258
$ws.Range("M32").Value = -8083.322
$ws.Range("N32").Value = -27749.258

# Hunk 11: sheet ARM, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2342.25
$ws.Range("J45").Value = 2928.4285
$ws.Range("L45").Value = 2928.4285
$ws.Range("N45").Value = -3682.4285

# Hunk 12: sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1793
$ws.Range("I61").Value = 1822.5333
$ws.Range("K61").Value = 1822.5333
$ws.Range("M61").Value = -1610.5333

# Hunk 13: sheet ARM, row 95
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 52374.75
$ws.Range("J95").Value = 52374.75
$ws.Range("L95").Value = 52374.75
$ws.Range("N95").Value = -57866.75

# Hunk 14: sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1793
$ws.Range("I136").Value = 1822.5333
$ws.Range("K136").Value = 5467.5999
$ws.Range("M136").Value = -2917.5999

# Hunk 15: sheet BSM, row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4205.6665
$ws.Range("I105").Value = 3484.2222
$ws.Range("K105").Value = 3484.2222
$ws.Range("M105").Value = -1737.2222

# Hunk 16: sheet BSM, row 115
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H115").Value = 64999
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

# Hunk 17: sheet CRP, row 25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 2000
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

# Hunk 18: sheet CRP, row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 52216.332
$ws.Range("J62").Value = 73328.336
$ws.Range("L62").Value = 73328.336
$ws.Range("N62").Value = -74576.336

# Hunk 19: sheet CRP, row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 52216.332
$ws.Range("J65").Value = 73328.336
$ws.Range("L65").Value = 366641.68
$ws.Range("N65").Value = -372881.68

# Hunk 20: sheet CRP, row 74
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 85499.5
$ws.Range("J74").Value = 85499.5
$ws.Range("L74").Value = 85499.5
$ws.Range("N74").Value = -87247.5

# Hunk 21: sheet CRP, row 77
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 85499.5
$ws.Range("J77").Value = 85499.5
$ws.Range("L77").Value = 256498.5
$ws.Range("N77").Value = -265234.5

# Hunk 22: sheet CRP, row 103
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 14743.889
$ws.Range("I103").Value = 12836.875
$ws.Range("K103").Value = 12836.875
$ws.Range("M103").Value = -11664.875

# Hunk 23: sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1586.3334
$ws.Range("I132").Value = 1586.3334
$ws.Range("K132").Value = 4759.0002
$ws.Range("M132").Value = -2229.0002

# Hunk 24: sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2491.1052
$ws.Range("I134").Value = 2013.75
$ws.Range("J134").Value = 5037
$ws.Range("K134").Value = 6041.25
$ws.Range("L134").Value = 15111
$ws.Range("M134").Value = -3506.25
$ws.Range("N134").Value = -20181

# Hunk 25: sheet CUL, row 14
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 2246.6667
$ws.Range("I14").Value = 2246.6667
$ws.Range("K14").Value = 6740.000100000001
$ws.Range("M14").Value = -6567.000100000001

# Hunk 26: sheet CUL, row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3557.1428
$ws.Range("J39").Value = 4975
$ws.Range("L39").Value = 14925
$ws.Range("N39").Value = -15513

# Hunk 27: sheet CUL, row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1169.7
$ws.Range("J55").Value = 1242.8572
$ws.Range("L55").Value = 3728.5716
$ws.Range("N55").Value = -4082.5716

# Hunk 28: sheet CUL, row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 11000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 11000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 99000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -104060

# Hunk 29: sheet CUL, row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 2389.1428
$ws.Range("I134").Value = 620.8333
$ws.Range("K134").Value = 1862.4999
$ws.Range("M134").Value = 3207.5001

# Hunk 30: sheet GSM, row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 625.6667
$ws.Range("I2").Value = 156.77777
$ws.Range("J2").Value = 2032.3334
$ws.Range("K2").Value = 156.77777
$ws.Range("L2").Value = 2032.3334
$ws.Range("M2").Value = -43.77777
$ws.Range("N2").Value = -2258.3334

# Hunk 31: sheet GSM, row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 30004
$ws.Range("I5").Value = 30004
$ws.Range("K5").Value = 30004
$ws.Range("M5").Value = -29892

# Hunk 32: sheet GSM, row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6641.364
$ws.Range("I70").Value = 4686.6665
$ws.Range("K70").Value = 4686.6665
$ws.Range("M70").Value = -4416.6665

# Hunk 33: sheet GSM, row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6641.364
$ws.Range("I73").Value = 4686.6665
$ws.Range("K73").Value = 4686.6665
$ws.Range("M73").Value = -3750.6665

# Hunk 34: sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9088.333000000001
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

# Hunk 35: sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 9088.333000000001
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

# Hunk 36: sheet LTW, row 2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5052499.5
$ws.Range("I2").Value = 10005000
$ws.Range("K2").Value = 10005000
$ws.Range("M2").Value = -10004888

# Hunk 37: sheet LTW, row 80
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 30999.666
$ws.Range("I80").Value = 25000
$ws.Range("J80").Value = 33999.5
$ws.Range("K80").Value = 25000
$ws.Range("L80").Value = 33999.5
$ws.Range("M80").Value = -23877
$ws.Range("N80").Value = -36245.5

# Hunk 38: sheet LTW, row 83
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H83").Value = 30999.666
$ws.Range("I83").Value = 25000
$ws.Range("J83").Value = 33999.5
$ws.Range("K83").Value = 75000
$ws.Range("L83").Value = 101998.5
$ws.Range("M83").Value = -69384
$ws.Range("N83").Value = -113230.5

# Hunk 39: sheet WVR, row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 20000
$ws.Range("J54").Value = 20000
$ws.Range("L54").Value = 20000
$ws.Range("N54").Value = -21040

# Hunk 40: sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 58508.223
$ws.Range("I136").Value = 1688.3334
$ws.Range("K136").Value = 5065.0002
$ws.Range("M136").Value = -2515.0002
